# Demo Robot with attachment validation
# - update the outlook account used by the robot (strEmailAccount row)
# - leave the "Description" column (C4) untouched
# - restore the selection to where the user left off (B18:B19) before saving

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "Sean.Crotty@defra.gov.uk"

[void]$ws.Range("B18:B19").Select()
